# Week 15 logged + Week 16 simulated: append this week's per-play/per-drive
# logs to the running season lists (YDS, ST) and update the cumulative
# season totals (OFF, DEF, ST, TURNS, PEN).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# YDS sheet: append newly logged play-by-play yardage figures to the
# season-long space-separated lists.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = $ws.Range("B2").Value2 + " 2 3 6 2 0 3 4 2 0 8 3 6 3 4 4 3 6 3 -1 0 2 5 4 2 1 0"
$ws.Range("C2").Value = $ws.Range("C2").Value2 + " 10 5 6 3 7 22 7 7 1 7 6 12 11 33 2 6 18 3 43"
$ws.Range("B3").Value = $ws.Range("B3").Value2 + " 3 10 3 0 7 12 16 0 4 1 2 -1 8 8 5 1 1 9 3 -1 2 5 -1 10"
$ws.Range("C3").Value = $ws.Range("C3").Value2 + " 30 11 12 21 9 8 3 6 5 10 11 2 11 7 5 5 7 5 5 11"

# ---------------------------------------------------------------------
# OFF sheet: updated cumulative season totals for Home (row 2) / Road (row 3)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("OFF")

$ws.Range("B2").Value = 7
$ws.Range("C2").Value = 142
$ws.Range("F2").Value = 71
$ws.Range("G2").Value = 34
$ws.Range("J2").Value = 27
$ws.Range("N2").Value = 14

$ws.Range("B3").Value = 7
$ws.Range("C3").Value = 127
$ws.Range("E3").Value = 27
$ws.Range("F3").Value = 78
$ws.Range("G3").Value = 27
$ws.Range("H3").Value = 36
$ws.Range("I3").Value = 54
$ws.Range("J3").Value = 55
$ws.Range("L3").Value = 221
$ws.Range("M3").Value = 137
$ws.Range("Q3").Value = 408

# ---------------------------------------------------------------------
# DEF sheet: updated cumulative season totals for Home (row 2) / Road (row 3)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("DEF")

$ws.Range("C2").Value = 198
$ws.Range("D2").Value = 10
$ws.Range("E2").Value = 11
$ws.Range("F2").Value = 46
$ws.Range("G2").Value = 55
$ws.Range("I2").Value = 5
$ws.Range("J2").Value = 27
$ws.Range("N2").Value = 18
$ws.Range("O2").Value = 17
$ws.Range("P2").Value = 11

$ws.Range("B3").Value = 10
$ws.Range("C3").Value = 131
$ws.Range("E3").Value = 31
$ws.Range("F3").Value = 93
$ws.Range("G3").Value = 22
$ws.Range("H3").Value = 23
$ws.Range("I3").Value = 43
$ws.Range("J3").Value = 43
$ws.Range("L3").Value = 230
$ws.Range("M3").Value = 153
$ws.Range("Q3").Value = 449

# ---------------------------------------------------------------------
# ST sheet: updated cumulative season totals (row 2/3) and per-kick
# season-long lists (B4:B6, D3:D5)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ST")

$ws.Range("B2").Value = 49
$ws.Range("D2").Value = 74
$ws.Range("F2").Value = 177
$ws.Range("G2").Value = 159
$ws.Range("J2").Value = 74
$ws.Range("K2").Value = 71
$ws.Range("N2").Value = 28
$ws.Range("O2").Value = 17

$ws.Range("B3").Value = 27

$ws.Range("D3").Value = $ws.Range("D3").Value2 + " 49 55 51 38"
$ws.Range("B4").Value = $ws.Range("B4").Value2 + " 68 67"
$ws.Range("D4").Value = $ws.Range("D4").Value2 + " -2 48 3 0"
$ws.Range("B5").Value = $ws.Range("B5").Value2 + " 20 16"
$ws.Range("D5").Value = $ws.Range("D5").Value2 + " 16 3 0 0 14"
$ws.Range("B6").Value = $ws.Range("B6").Value2 + " 98 0 23"

# ---------------------------------------------------------------------
# TURNS sheet: updated cumulative season totals for Road (row 3)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("TURNS")

$ws.Range("B3").Value = 11
$ws.Range("E3").Value = 6

# ---------------------------------------------------------------------
# PEN sheet: updated cumulative season totals for Holding (row 3)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("PEN")

$ws.Range("B3").Value = 22
$ws.Range("D3").Value = 5
